# Adds newly-measured embryo/endosperm weight rows for three more species
# (BAER, PEPU, GRBU), labels the data with a header, and re-sorts the whole
# table alphabetically by species code - matching the commit
# "adding file with embryo+endosperm weights determined from my species".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append the new species' raw measurements below the existing 11 rows
#    of data (rows 13-20), typed in the order BAER, then PEPU, then GRBU -
#    this is the order the new shared strings need to be introduced in.
$ws.Range("A13").Value = "BAER"
$ws.Range("B13").Value = 23.463999999999999
$ws.Range("A14").Value = "BAER"
$ws.Range("B14").Value = 24.664999999999999

$ws.Range("A15").Value = "PEPU"
$ws.Range("B15").Value = 2.2200000000000002
$ws.Range("A16").Value = "PEPU"
$ws.Range("B16").Value = 2.028
$ws.Range("A17").Value = "PEPU"
$ws.Range("B17").Value = 2.3820000000000001

$ws.Range("A18").Value = "GRBU"
$ws.Range("B18").Value = 26.265000000000001
$ws.Range("A19").Value = "GRBU"
$ws.Range("B19").Value = 28.085999999999999
$ws.Range("A20").Value = "GRBU"
$ws.Range("B20").Value = 25.745000000000001

# 2) Column A had no header before - label it now that the sheet holds
#    more than one species' worth of rows.
$ws.Range("A1").Value = "species"

# 3) Sort all the data (A2:B20) ascending by species code, with a header
#    row above it, the same as Data > Sort in the Excel UI.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A20"))
$ws.Sort.SetRange($ws.Range("A1:B20"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# 4) Column B now holds a long header ("endosperm_embryo_size") - widen it
#    to fit, as Excel does automatically when you widen a column to show it.
$ws.Columns.Item(2).AutoFit() | Out-Null

# 5) Leave the selection where it ended up after entering/reviewing the
#    newly-added GRBU rows.
$ws.Range("B18:B20").Select() | Out-Null
